# Rewrites the input-list generator output: regenerates the per-trial
# rows of the kitchens categorization stimulus list (sheet1) so that all
# subjects receive the same selection of images, adding one additional
# trial (row 28) and shifting/replacing the target/distractor assignments,
# stimuli file names and associated conceptual/perceptual/typicality stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 1
    $row[0,5] = 28
    $row[0,6] = "kitchens"
    $row[0,7] = "bedrooms"
    $row[0,8] = "distractor"
    $row[0,9] = $null
    $row[0,10] = "f"
    $row[0,11] = "stimuli/img_4wq98.png"
    $row[0,12] = 78.48387096774194
    $row[0,13] = 58.12903225806452
    $row[0,14] = 68.30645161290323
    $row[0,15] = 31
    $row[0,16] = 7
    $row[0,17] = 7
    $row[0,18] = 7
    $ws.Range("A2:S2").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 2
    $row[0,5] = 29
    $row[0,6] = "kitchens"
    $row[0,7] = "living_rooms"
    $row[0,8] = "distractor"
    $row[0,9] = $null
    $row[0,10] = "f"
    $row[0,11] = "stimuli/img_eh0no.png"
    $row[0,12] = 53.66666666666666
    $row[0,13] = 36.02564102564103
    $row[0,14] = 44.84615384615385
    $row[0,15] = 39
    $row[0,16] = 3
    $row[0,17] = 3
    $row[0,18] = 3
    $ws.Range("A3:S3").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 3
    $row[0,5] = 30
    $row[0,6] = "kitchens"
    $row[0,7] = "bedrooms"
    $row[0,8] = "distractor"
    $row[0,9] = $null
    $row[0,10] = "f"
    $row[0,11] = "stimuli/img_ca8kd.png"
    $row[0,12] = 92.05405405405405
    $row[0,13] = 73.02702702702703
    $row[0,14] = 82.54054054054055
    $row[0,15] = 37
    $row[0,16] = 10
    $row[0,17] = 10
    $row[0,18] = 10
    $ws.Range("A4:S4").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 4
    $row[0,5] = 31
    $row[0,6] = "kitchens"
    $row[0,7] = "living_rooms"
    $row[0,8] = "distractor"
    $row[0,9] = $null
    $row[0,10] = "f"
    $row[0,11] = "stimuli/img_x9w7o.png"
    $row[0,12] = 92.38888888888889
    $row[0,13] = 72.94444444444444
    $row[0,14] = 82.66666666666666
    $row[0,15] = 36
    $row[0,16] = 10
    $row[0,17] = 10
    $row[0,18] = 10
    $ws.Range("A5:S5").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 5
    $row[0,5] = 32
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_ye5sl.png"
    $row[0,12] = 53.2258064516129
    $row[0,13] = 34.45161290322581
    $row[0,14] = 43.83870967741936
    $row[0,15] = 31
    $row[0,16] = 2
    $row[0,17] = 2
    $row[0,18] = 2
    $ws.Range("A6:S6").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 6
    $row[0,5] = 33
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_es7o2.png"
    $row[0,12] = 52.48571428571429
    $row[0,13] = 27.54285714285714
    $row[0,14] = 40.01428571428572
    $row[0,15] = 35
    $row[0,16] = 2
    $row[0,17] = 2
    $row[0,18] = 2
    $ws.Range("A7:S7").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 7
    $row[0,5] = 34
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_p3hpc.png"
    $row[0,12] = 72.83333333333333
    $row[0,13] = 52.22222222222222
    $row[0,14] = 62.52777777777777
    $row[0,15] = 36
    $row[0,16] = 6
    $row[0,17] = 6
    $row[0,18] = 6
    $ws.Range("A8:S8").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 8
    $row[0,5] = 35
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_wyl6z.png"
    $row[0,12] = 59.8235294117647
    $row[0,13] = 36.23529411764706
    $row[0,14] = 48.02941176470588
    $row[0,15] = 34
    $row[0,16] = 3
    $row[0,17] = 3
    $row[0,18] = 3
    $ws.Range("A9:S9").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 9
    $row[0,5] = 36
    $row[0,6] = "kitchens"
    $row[0,7] = "bedrooms"
    $row[0,8] = "distractor"
    $row[0,9] = $null
    $row[0,10] = "f"
    $row[0,11] = "stimuli/img_72fmj.png"
    $row[0,12] = 53.87179487179487
    $row[0,13] = 36.02564102564103
    $row[0,14] = 44.94871794871795
    $row[0,15] = 39
    $row[0,16] = 3
    $row[0,17] = 3
    $row[0,18] = 3
    $ws.Range("A10:S10").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 10
    $row[0,5] = 37
    $row[0,6] = "kitchens"
    $row[0,7] = "living_rooms"
    $row[0,8] = "distractor"
    $row[0,9] = $null
    $row[0,10] = "f"
    $row[0,11] = "stimuli/img_jpjeg.png"
    $row[0,12] = 90.90697674418605
    $row[0,13] = 74.3953488372093
    $row[0,14] = 82.65116279069767
    $row[0,15] = 43
    $row[0,16] = 10
    $row[0,17] = 10
    $row[0,18] = 10
    $ws.Range("A11:S11").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 11
    $row[0,5] = 38
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_60242.png"
    $row[0,12] = 78.33333333333333
    $row[0,13] = 57.57575757575758
    $row[0,14] = 67.95454545454545
    $row[0,15] = 33
    $row[0,16] = 7
    $row[0,17] = 7
    $row[0,18] = 7
    $ws.Range("A12:S12").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 12
    $row[0,5] = 39
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_t90e2.png"
    $row[0,12] = 83.0625
    $row[0,13] = 61.96875
    $row[0,14] = 72.515625
    $row[0,15] = 32
    $row[0,16] = 9
    $row[0,17] = 9
    $row[0,18] = 9
    $ws.Range("A13:S13").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 13
    $row[0,5] = 40
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_7wul8.png"
    $row[0,12] = 43.03030303030303
    $row[0,13] = 25.54545454545455
    $row[0,14] = 34.28787878787879
    $row[0,15] = 33
    $row[0,16] = 1
    $row[0,17] = 1
    $row[0,18] = 1
    $ws.Range("A14:S14").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 14
    $row[0,5] = 41
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_ce9vx.png"
    $row[0,12] = 75.9090909090909
    $row[0,13] = 57.12121212121212
    $row[0,14] = 66.51515151515152
    $row[0,15] = 33
    $row[0,16] = 7
    $row[0,17] = 7
    $row[0,18] = 7
    $ws.Range("A15:S15").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 15
    $row[0,5] = 42
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_cnyac.png"
    $row[0,12] = 69.1470588235294
    $row[0,13] = 47.8235294117647
    $row[0,14] = 58.48529411764706
    $row[0,15] = 34
    $row[0,16] = 5
    $row[0,17] = 5
    $row[0,18] = 5
    $ws.Range("A16:S16").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 16
    $row[0,5] = 43
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_a8wvq.png"
    $row[0,12] = 86.25925925925925
    $row[0,13] = 66.25925925925925
    $row[0,14] = 76.25925925925925
    $row[0,15] = 27
    $row[0,16] = 10
    $row[0,17] = 10
    $row[0,18] = 10
    $ws.Range("A17:S17").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 17
    $row[0,5] = 44
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_aplao.png"
    $row[0,12] = 64.0909090909091
    $row[0,13] = 40.75757575757576
    $row[0,14] = 52.42424242424242
    $row[0,15] = 33
    $row[0,16] = 3
    $row[0,17] = 3
    $row[0,18] = 3
    $ws.Range("A18:S18").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 18
    $row[0,5] = 45
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_nyv2b.png"
    $row[0,12] = 11.91176470588235
    $row[0,13] = 6.852941176470588
    $row[0,14] = 9.382352941176471
    $row[0,15] = 34
    $row[0,16] = 1
    $row[0,17] = 1
    $row[0,18] = 1
    $ws.Range("A19:S19").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 19
    $row[0,5] = 46
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_inqod.png"
    $row[0,12] = 70.84848484848484
    $row[0,13] = 50.63636363636363
    $row[0,14] = 60.74242424242424
    $row[0,15] = 33
    $row[0,16] = 5
    $row[0,17] = 5
    $row[0,18] = 5
    $ws.Range("A20:S20").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 20
    $row[0,5] = 47
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_yeh72.png"
    $row[0,12] = 68.66666666666667
    $row[0,13] = 45.21212121212121
    $row[0,14] = 56.93939393939394
    $row[0,15] = 33
    $row[0,16] = 4
    $row[0,17] = 4
    $row[0,18] = 4
    $ws.Range("A21:S21").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 21
    $row[0,5] = 48
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_cv6mf.png"
    $row[0,12] = 66.8
    $row[0,13] = 42.08
    $row[0,14] = 54.44
    $row[0,15] = 25
    $row[0,16] = 4
    $row[0,17] = 4
    $row[0,18] = 4
    $ws.Range("A22:S22").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 22
    $row[0,5] = 49
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_d8xbu.png"
    $row[0,12] = 91.36363636363636
    $row[0,13] = 73.18181818181819
    $row[0,14] = 82.27272727272728
    $row[0,15] = 33
    $row[0,16] = 10
    $row[0,17] = 10
    $row[0,18] = 10
    $ws.Range("A23:S23").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 23
    $row[0,5] = 50
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_iyxnj.png"
    $row[0,12] = 75.30555555555556
    $row[0,13] = 54.33333333333334
    $row[0,14] = 64.81944444444444
    $row[0,15] = 36
    $row[0,16] = 6
    $row[0,17] = 6
    $row[0,18] = 6
    $ws.Range("A24:S24").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 24
    $row[0,5] = 51
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_9mky8.png"
    $row[0,12] = 84.32352941176471
    $row[0,13] = 65.17647058823529
    $row[0,14] = 74.75
    $row[0,15] = 34
    $row[0,16] = 9
    $row[0,17] = 9
    $row[0,18] = 9
    $ws.Range("A25:S25").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 25
    $row[0,5] = 52
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_uwv6y.png"
    $row[0,12] = 78.88888888888889
    $row[0,13] = 59.30555555555556
    $row[0,14] = 69.09722222222223
    $row[0,15] = 36
    $row[0,16] = 8
    $row[0,17] = 8
    $row[0,18] = 8
    $ws.Range("A26:S26").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 26
    $row[0,5] = 53
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_eatdk.png"
    $row[0,12] = 81.40625
    $row[0,13] = 61.375
    $row[0,14] = 71.390625
    $row[0,15] = 32
    $row[0,16] = 8
    $row[0,17] = 8
    $row[0,18] = 8
    $ws.Range("A27:S27").Value = $row

    $row = New-Object 'object[,]' 1,19
    $row[0,0] = 1
    $row[0,1] = "categorization"
    $row[0,2] = 1
    $row[0,3] = 1
    $row[0,4] = 27
    $row[0,5] = 54
    $row[0,6] = "kitchens"
    $row[0,7] = "kitchens"
    $row[0,8] = "target"
    $row[0,9] = $null
    $row[0,10] = "j"
    $row[0,11] = "stimuli/img_r77yy.png"
    $row[0,12] = 84.125
    $row[0,13] = 65.375
    $row[0,14] = 74.75
    $row[0,15] = 32
    $row[0,16] = 9
    $row[0,17] = 9
    $row[0,18] = 9
    $ws.Range("A28:S28").Value = $row
